$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.548.93"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "3.442.17"
$ws.Range("E3").Value = "  -2.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.31"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.65"
$ws.Range("E6").Value = "  -6.54%  "
$ws.Range("D7").Value = "3.442.23"
$ws.Range("E7").Value = "  -2.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("E10").Value = "  -5.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.123"
$ws.Range("E11").Value = "  -8.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.377"
$ws.Range("E12").Value = "  -7.46%  "
$ws.Range("D13").Value = "4.022.50"
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("E14").Value = "  -10.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.61"
$ws.Range("E15").Value = "  -8.74%  "
$ws.Range("D16").Value = "3.476.42"
$ws.Range("E16").Value = "  -1.74%  "
$ws.Range("D17").Value = "65.526.24"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.92"
$ws.Range("E19").Value = "  -10.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.83"
$ws.Range("E20").Value = "  -6.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.68"
$ws.Range("E21").Value = "  -6.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "394.26"
$ws.Range("E22").Value = "  -5.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.549"
$ws.Range("E23").Value = "  -8.92%  "
$ws.Range("E24").Value = "  -5.40%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "3.583.83"
$ws.Range("E26").Value = "  -2.82%  "
$ws.Range("E27").Value = "  -9.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.22"
$ws.Range("E29").Value = "  -8.18%  "
$ws.Range("E30").Value = "  -8.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.20"
$ws.Range("E31").Value = "  -10.52%  "
$ws.Range("D32").Value = "3.447.79"
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -5.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.01"
$ws.Range("E35").Value = "  -6.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "171.24"
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("E37").Value = "  -8.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.21"
$ws.Range("E38").Value = "  -10.92%  "
$ws.Range("E39").Value = "  -6.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.85"
$ws.Range("E40").Value = "  -9.68%  "
$ws.Range("E41").Value = "  -7.03%  "
$ws.Range("E42").Value = "  -4.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.59"
$ws.Range("E43").Value = "  -4.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.43"
$ws.Range("E45").Value = "  -13.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.62"
$ws.Range("E46").Value = "  -10.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.10"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.66"
$ws.Range("E48").Value = "  -2.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.55"
$ws.Range("E49").Value = "  -7.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.10"
$ws.Range("E50").Value = "  -14.04%  "
$ws.Range("D51").Value = "2.204.06"
$ws.Range("E51").Value = "  -6.81%  "
